$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column K: header "Revenues", styled like the other header cells ---
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K1").Value = "Revenues"

# --- Row 2 updates ---
# Plain text updates (values Excel wouldn't reinterpret as numbers)
$ws.Range("C2").Value = "+0.19(0.26%) 1D"
$ws.Range("D2").Value = "$309.11B"

# Values that look numeric to Excel's auto-detection need to be forced to stay
# as literal text (matching the source inlineStr cells), without leaving any
# residual number-format/style change on the cell.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "$72.06"
$ws.Range("B2").ClearFormats()

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "29.1"
$ws.Range("F2").ClearFormats()

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2.84%"
$ws.Range("I2").ClearFormats()

# --- New data cell K2 ---
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "10,980"
$ws.Range("K2").ClearFormats()
